$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This revision touches no visible text. It only refreshes two kinds of
# internal marker ids that M2Doc stamps on the generated document each time
# it is (re)produced:
#
#   1. The w:rsidR stamped on the five runs that make up the
#      "REF bookmark1 \h" field (begin / instrText / separate / result / end)
#      in the "Test link before bookmark : ..." paragraph.
#   2. The w:id pair on the bookmarkStart/bookmarkEnd that wrap the
#      "Invalid block: ..." bookmark1 target text.
#
# Neither value is exposed as a settable property on the regular Word object
# model, so the field-run rsid is refreshed by re-submitting that paragraph's
# markup through Range.InsertXML (everything else in the paragraph is kept
# byte-for-byte identical), and the bookmark id is refreshed the idiomatic
# way: delete the bookmark and re-Add it over the same range, which mints a
# fresh id for it.
# ---------------------------------------------------------------------------

$newRunId = "196BDFDDB4ED4956AEFC1959D8212C2E"

# --- 2. Refresh bookmark1's id (delete + re-add over the same range) -------
$bm = $d.Bookmarks.Item("bookmark1")
$bmRange = $bm.Range.Duplicate
$bm.Delete()
$d.Bookmarks.Add("bookmark1", $bmRange)

# --- 1. Refresh the REF field runs' rsidR -----------------------------------
$p2 = $d.Paragraphs(2)
$p2xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidP="009168BC" w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B" w:rsidRPr="00FF681D"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00FF681D"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Test link before bookmark : </w:t></w:r><w:r w:rsidR="$newRunId"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="$newRunId"><w:instrText xml:space="preserve"> REF bookmark1 \h </w:instrText></w:r><w:r w:rsidR="$newRunId"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="$newRunId"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r><w:r w:rsidR="$newRunId"><w:fldChar w:fldCharType="end"/></w:r></w:p>
"@
$p2.Range.InsertXML($p2xml)

Write-Output "rsidR -> $newRunId on the REF bookmark1 field runs; bookmark1 id refreshed"
